$wb = $excel.ActiveWorkbook
$members = $wb | Get-Member
Write-Host $members
